# The presentation's design was switched from the custom "Integral" theme
# to the default PowerPoint "Office Theme". Apply this by updating the
# slide master's theme: its name and its 12 theme colors (the font scheme
# and format scheme are identical between the two themes, so only the
# color scheme actually changes).

$p = $ppt.ActivePresentation

function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$theme = $p.SlideMaster.Theme
$theme.Name = "Office Theme"

$tcs = $theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = HexToRgbInt($officeColors[$i - 1])
}

# The notes master previously carried the "Office Theme" colors while the
# slide master carried "Integral"; after the swap the notes master should
# reflect the "Integral" palette instead.
$integralColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "455F51",  # 3  dk2
    "E3DED1",  # 4  lt2
    "99CB38",  # 5  accent1
    "63A537",  # 6  accent2
    "E6D024",  # 7  accent3
    "CC9700",  # 8  accent4
    "4EB3CF",  # 9  accent5
    "378DA6",  # 10 accent6
    "6B9F25",  # 11 hlink
    "B26B02"   # 12 folHlink
)

try {
    $nmTheme = $p.NotesMaster.Theme
    $nmTheme.Name = "Integral"
    $nmTcs = $nmTheme.ThemeColorScheme
    for ($i = 1; $i -le $integralColors.Count; $i++) {
        $nmTcs.Item($i).RGB = HexToRgbInt($integralColors[$i - 1])
    }
} catch {
    # Notes master theme may not be independently addressable in all hosts.
}
